# fix(publipostage): Correct status name
#
# - "bleu" -> "noir" (statut_label column)
# - statut_name labels reworded from "... posté(e)(s) ..." to "... postés ou publiés ..."
#
# Longer/more specific strings are replaced before the shorter strings they
# contain, so the substring "résultat et / ou publication posté" doesn't
# get partially rewritten before the "dans les 36/12 mois" variants are
# handled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("bleu", "noir")

$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés")
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
